# Apply crypto price/volume updates scraped on Thu Mar 16 11:11:45 UTC 2023
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell([string]$addr, [string]$val) {
    $rng = $ws.Range($addr)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $origStyle
}

Set-TextCell "D2" "24.945.78"
Set-TextCell "E2" "  +1.19%  "

Set-TextCell "D3" "1.662.87"
Set-TextCell "E3" "  -1.33%  "

Set-TextCell "D4" "0.9997"
Set-TextCell "E4" "  -0.86%  "

Set-TextCell "D5" "317.36"
Set-TextCell "E5" "  +2.93%  "

Set-TextCell "D6" "0.9986"
Set-TextCell "E6" "  -0.38%  "

Set-TextCell "D7" "0.3636"
Set-TextCell "E7" "  -1.70%  "

Set-TextCell "D8" "47.15"
Set-TextCell "E8" "  -2.48%  "

Set-TextCell "D9" "0.3278"
Set-TextCell "E9" "  -2.76%  "

Set-TextCell "D10" "1.142"
Set-TextCell "E10" "  -3.68%  "

Set-TextCell "D11" "0.07078"
Set-TextCell "E11" "  -3.49%  "

Set-TextCell "D12" "0.9995"
Set-TextCell "E12" "  -0.47%  "

Set-TextCell "D13" "6.061"
Set-TextCell "E13" "  -2.23%  "

Set-TextCell "D14" "19.66"
Set-TextCell "E14" "  -4.29%  "

Set-TextCell "D15" "1.663.74"
Set-TextCell "E15" "  -1.46%  "

Set-TextCell "D16" "6.629"
Set-TextCell "E16" "  -3.32%  "

Set-TextCell "D17" "0.00001051"
Set-TextCell "E17" "  -4.63%  "

Set-TextCell "D18" "0.06627"
Set-TextCell "E18" "  -0.26%  "

Set-TextCell "D19" "0.9973"
Set-TextCell "E19" "  -0.43%  "

Set-TextCell "D20" "79.55"
Set-TextCell "E20" "  -2.88%  "

Set-TextCell "D21" "5.917"
Set-TextCell "E21" "  -4.87%  "

Set-TextCell "D22" "15.78"
Set-TextCell "E22" "  -6.54%  "

Set-TextCell "D23" "12.57"
Set-TextCell "E23" "  -1.51%  "

Set-TextCell "D24" "24.895.73"
Set-TextCell "E24" "  +0.94%  "

Set-TextCell "D25" "2.433"
Set-TextCell "E25" "  +0.02%  "

Set-TextCell "D26" "2.405"
Set-TextCell "E26" "  -11.00%  "

Set-TextCell "D27" "148.47"
Set-TextCell "E27" "  +0.75%  "

Set-TextCell "D28" "18.67"
Set-TextCell "E28" "  -6.13%  "

Set-TextCell "B29" "ImmutableX"
Set-TextCell "C29" "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextCell "D29" "1.232"
Set-TextCell "E29" "  +2.15%  "

Set-TextCell "B30" "WrappedliquidstakedEther2.0"
Set-TextCell "C30" "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
Set-TextCell "D30" "1.850.78"
Set-TextCell "E30" "  -1.26%  "

Set-TextCell "D31" "125.75"
Set-TextCell "E31" "  -3.63%  "

Set-TextCell "D32" "4.138"
Set-TextCell "E32" "  -1.17%  "

Set-TextCell "D33" "5.849"
Set-TextCell "E33" "  -10.37%  "

Set-TextCell "D34" "0.08455"
Set-TextCell "E34" "  -1.67%  "

Set-TextCell "D35" "1.674"
Set-TextCell "E35" "  -3.20%  "

Set-TextCell "E36" "  -7.71%  "

Set-TextCell "D37" "1.283"
Set-TextCell "E37" "  +2.83%  "

Set-TextCell "D38" "5.216"
Set-TextCell "E38" "  -3.65%  "

Set-TextCell "D39" "0.02261"
Set-TextCell "E39" "  -3.76%  "

Set-TextCell "D40" "0.06076"
Set-TextCell "E40" "  -6.41%  "

Set-TextCell "B41" "FraxShare"
Set-TextCell "C41" "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextCell "D41" "8.329"
Set-TextCell "E41" "  -6.14%  "

Set-TextCell "B42" "Algorand"
Set-TextCell "C42" "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextCell "D42" "0.2072"
Set-TextCell "E42" "  -4.58%  "

Set-TextCell "D43" "0.9973"
Set-TextCell "E43" "  -0.43%  "

Set-TextCell "D44" "0.5950"
Set-TextCell "E44" "  -5.13%  "

Set-TextCell "D45" "3.820"
Set-TextCell "E45" "  +1.09%  "

Set-TextCell "D46" "12.80"
Set-TextCell "E46" "  -4.75%  "

Set-TextCell "D47" "0.5645"
Set-TextCell "E47" "  -5.14%  "

Set-TextCell "D48" "125.94"
Set-TextCell "E48" "  -0.52%  "

Set-TextCell "D49" "1.962"
Set-TextCell "E49" "  -4.78%  "

Set-TextCell "D50" "0.07027"
Set-TextCell "E50" "  -2.14%  "

Set-TextCell "D51" "1.196"
Set-TextCell "E51" "  -0.05%  "
